$d = $word.ActiveDocument

# Helper: replace a single character immediately following the first
# occurrence of $anchor (searched starting at char position $searchStart),
# used for the MERGEFIELD-produced day-of-month digit which is its own
# isolated run wedged between field-character runs. A plain whole-document
# Find/Replace on the bare digit is unsafe because the digit also occurs
# inside unrelated numeric runs (e.g. "53.458"). Returns the position right
# after the replaced character, so callers can keep searching forward for
# subsequent occurrences without re-matching the one just fixed.
function Replace-FieldDigit($doc, $anchor, $newDigit, $searchStart) {
    $scope = $doc.Range($searchStart, $doc.Content.End)
    $found = $scope.Find.Execute($anchor, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $digitStart = $scope.Start
    $digitRange = $doc.Range($digitStart, $digitStart + 1)
    $digitRange.Text = $newDigit
    return $digitStart + 1
}

# 1. Salutation MS -> DR (both copies of the certificate)
$d.Content.Find.Execute("MS", $true, $false, $false, $false, $false, $true, 1, $false, "DR", 2)

# 2. Employee full name
$d.Content.Find.Execute("FELICITAS M. SUMAGUI", $true, $false, $false, $false, $false, $true, 1, $false, "LIZA FE P. CAPUPUS", 2)

# 3. Position / designation
$d.Content.Find.Execute("Casual Employee", $true, $false, $false, $false, $false, $true, 1, $false, "City Health Officer II", 2)

# 4. Office
$d.Content.Find.Execute("City Social Welfare Development Office", $true, $false, $false, $false, $false, $true, 1, $false, "City Health Office", 2)

# 5. "as of" date
$d.Content.Find.Execute("March 22, 2023", $true, $false, $false, $false, $false, $true, 1, $false, "March 18, 2024", 2)

# 6. Vacation leave credits
$d.Content.Find.Execute("  53.458", $true, $false, $false, $false, $false, $true, 1, $false, "  90.542", 2)

# 7. Sick leave credits
$d.Content.Find.Execute("  67.458", $true, $false, $false, $false, $false, $true, 1, $false, " 101.542", 2)

# 8. Total leave credits
$d.Content.Find.Execute(" 120.916", $true, $false, $false, $false, $false, $true, 1, $false, " 192.084", 2)

# 9. Salutation Ms -> Dr (requester line, match case so "MS" above is untouched)
$d.Content.Find.Execute("Ms", $true, $false, $false, $false, $false, $true, 1, $false, "Dr", 2)

# 10. Requester surname
$d.Content.Find.Execute("Sumagui", $true, $false, $false, $false, $false, $true, 1, $false, "Capupus", 2)

# 11. "Issued this N(th|nd) day of" -> day-of-month digit, 5 -> 7 (both copies).
#     Scoped/targeted so it cannot touch the "5" digits inside the leave-credit
#     numbers above.
$pos = 0
$pos = Replace-FieldDigit $d "5th" "7" $pos
$pos = Replace-FieldDigit $d "5nd" "7" $pos

# 12. Month name July -> August
$d.Content.Find.Execute("July", $true, $false, $false, $false, $false, $true, 1, $false, "August", 2)
